$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "24/10/2025"
$ws.Range("B10").Value = "Paris FC"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = "Nantes"
$ws.Range("F10").Value = "W"
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0.92
$ws.Range("L10").Value = 1.16
$ws.Range("M10").Value = 11
$ws.Range("N10").Value = 14
$ws.Range("O10").Value = 3
$ws.Range("P10").Value = 3
